# Applies the "Added some client questions" edit:
#  - Splits the old "Judge Dashboard Flow" answer paragraph into a
#    "Mentor Evaluation Sheet" sub-section and adds two new sub-sections
#    ("Reports" and "Application Review") with their own questions.
#  - Merges a couple of runs that had been split mid-sentence.
#  - Moves the (cosmetic) lastRenderedPageBreak marker from the "Pitch
#    Evaluation" heading run onto the "Would you want to show ..." run,
#    matching the reflowed pagination.
#  - Marks the DefaultParagraphFont style as semiHidden.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1. Expand the "Describe or send us the judges' mentor evaluation
#    sheet..." paragraph into three sub-sections.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Describe or send us the judges' mentor evaluation sheet.*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $p = $d.Paragraphs.Item($target)
    $xml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading2"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Mentor Evaluation Sheet</w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:left="1440"/></w:pPr>' +
           '<w:r><w:t xml:space="preserve">Describe or send us the judges'' mentor evaluation sheet. </w:t></w:r>' +
           '<w:proofErr w:type="gramStart"/><w:r><w:t>We''re</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
           '<w:r><w:t xml:space="preserve"> not sure how a judge relates to a mentor, and therefore by extension, how a judge should evaluate a mentor. We </w:t></w:r>' +
           '<w:proofErr w:type="gramStart"/><w:r><w:t>can''t</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
           '<w:r><w:t xml:space="preserve"> build a UI for this without this understanding.</w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading2"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Reports</w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:left="1440"/></w:pPr><w:r><w:t xml:space="preserve">What sorts of reports will judges run?  </w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading2"/><w:ind w:left="720"/></w:pPr><w:r><w:t>Application Review</w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:left="1440"/></w:pPr><w:r><w:t>Are judges involved in application review?  If so, we can add it to the dashboard:</w:t></w:r></w:p>' +
           '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:left="2160"/></w:pPr><w:r><w:t xml:space="preserve"> What does the application review form look like if so?  Perhaps it is one of the documents we already have.</w:t></w:r></w:p>'
    $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 2. Merge the "1. " / "By Applicant - where in the process map " runs.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("1. By Applicant - where in the process map ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1. By Applicant - where in the process map ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Merge the "Would you want to show ... alongside ... the grading"
#    runs and move the lastRenderedPageBreak marker onto this run.
# ---------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Would you want to show the application on this page*") {
        $target2 = $i
        break
    }
}

if ($target2 -ne $null) {
    $p2 = $d.Paragraphs.Item($target2)
    $xml2 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:left="720"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>' +
            '<w:r><w:lastRenderedPageBreak/><w:t>Would you want to show the application on this page alongside the grading</w:t></w:r>' +
            '<w:r><w:t>?</w:t></w:r>' +
            '<w:r><w:br/></w:r>' +
            '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Compilation note: the grading for what?  Are we talking about the judging process?</w:t></w:r></w:p>'
    $p2.Range.InsertXML($xml2)
}

# ---------------------------------------------------------------------
# 4. Remove the lastRenderedPageBreak marker from "Pitch Evaluation".
# ---------------------------------------------------------------------
$target3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Pitch Evaluation`r") {
        $target3 = $i
        break
    }
}

if ($target3 -ne $null) {
    $p3 = $d.Paragraphs.Item($target3)
    $xml3 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Pitch Evaluation</w:t></w:r></w:p>'
    $p3.Range.InsertXML($xml3)
}

# ---------------------------------------------------------------------
# 5. Merge the "What else could be pulled up for managers ..." runs.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("What else could be pulled up for managers look at when the page is loaded?", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "What else could be pulled up for managers look at when the page is loaded?", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Mark the DefaultParagraphFont style as semiHidden.
# ---------------------------------------------------------------------
$style = $d.Styles.Item("Default Paragraph Font")
$style.SemiHidden = $true

Write-Output "done"
